# Proyecto Final Desarrollo ABAP con COPILOT.xlsx
# Fill in "Asistido Copilot" (column D) effort hours for Bloque 2: Desarrollo CORE
# rows 15-19, extend the block's D-total formula to include the newly-added
# row 19, and set the Responsable (column E) for the two rows that didn't
# have it yet. Downstream totals (D38/D39/D40) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = 2.5
$ws.Range("D16").Value = 2.5
$ws.Range("D17").Value = 2.5

$ws.Range("D18").Value = 2.5
$ws.Range("E18").Value = "Jhonatan Hidalgo"

$ws.Range("D19").Value = 4
$ws.Range("E19").Value = "Jhonatan Hidalgo"

# Bloque total now spans through the newly-populated row 19.
$ws.Range("D20").Formula = "=SUM(D15:D19)"

# Restore the last active-cell selection recorded in the sheet.
$ws.Range("G20").Select() | Out-Null
